$wb = $excel.ActiveWorkbook

# "Operand 1" sheet: C3 value 2 -> 3
$ws1 = $wb.Worksheets.Item("Operand 1")
$ws1.Range("C3").Value = 3

# Author switched the active/selected tab from "Explicit number ranges"
# to "Operand 2" (activeTab 2 -> 1, tabSelected moves sheets).
$ws2 = $wb.Worksheets.Item("Operand 2")
$ws2.Activate()
